$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 259.09525
$ws.Range("I28").Value = 153.125
$ws.Range("J28").Value = 598.2
$ws.Range("K28").Value = 153.125
$ws.Range("L28").Value = 598.2
$ws.Range("M28").Value = 331.875
$ws.Range("N28").Value = -1568.2

$ws.Range("H61").Value = 4178.2144
$ws.Range("I61").Value = 1193.5
$ws.Range("K61").Value = 3580.5
$ws.Range("M61").Value = -3408.5

$ws.Range("H94").Value = 1152
$ws.Range("I94").Value = 1152
$ws.Range("K94").Value = 1152
$ws.Range("M94").Value = -701

$ws.Range("H113").Value = 4075.0833
$ws.Range("I113").Value = 4150.722
$ws.Range("K113").Value = 4150.722
$ws.Range("M113").Value = -896.7219999999998

$ws.Range("H132").Value = 3349.242
$ws.Range("I132").Value = 2999.4255
$ws.Range("K132").Value = 8998.2765
$ws.Range("M132").Value = -6468.2765

$ws.Range("H137").Value = 9535.513000000001
$ws.Range("I137").Value = 5002.48
$ws.Range("K137").Value = 15007.44
$ws.Range("M137").Value = -12457.44

$ws.Range("H138").Value = 3034.7273
$ws.Range("I138").Value = 1988.2858
$ws.Range("J138").Value = 3316.4614
$ws.Range("K138").Value = 5964.857400000001
$ws.Range("L138").Value = 9949.3842
$ws.Range("M138").Value = -824.8574000000008
$ws.Range("N138").Value = -20229.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5439.4287
$ws.Range("I2").Value = 4035.1738
$ws.Range("K2").Value = 4035.1738
$ws.Range("M2").Value = -3922.1738

$ws.Range("H32").Value = 2684.5178
$ws.Range("I32").Value = 1510.1063
$ws.Range("K32").Value = 1510.1063
$ws.Range("M32").Value = -1223.1063

$ws.Range("H45").Value = 7945.4346
$ws.Range("I45").Value = 8795.684999999999
$ws.Range("K45").Value = 8795.684999999999
$ws.Range("M45").Value = -8418.684999999999

$ws.Range("H61").Value = 11889.917
$ws.Range("I61").Value = 8085.875
$ws.Range("J61").Value = 19498
$ws.Range("K61").Value = 8085.875
$ws.Range("L61").Value = 19498
$ws.Range("M61").Value = -7873.875
$ws.Range("N61").Value = -19922

$ws.Range("H74").Value = 3598.4375
$ws.Range("I74").Value = 1212
$ws.Range("K74").Value = 1212
$ws.Range("M74").Value = -338

$ws.Range("H77").Value = 3598.4375
$ws.Range("I77").Value = 1212
$ws.Range("K77").Value = 6060
$ws.Range("M77").Value = -1692

$ws.Range("H110").Value = 2133.8
$ws.Range("I110").Value = 2016.5294
$ws.Range("K110").Value = 2016.5294
$ws.Range("M110").Value = 28.4706000000001

$ws.Range("H116").Value = 5439.4287
$ws.Range("I116").Value = 4035.1738
$ws.Range("K116").Value = 4035.1738
$ws.Range("M116").Value = -1741.1738

$ws.Range("H132").Value = 9284.804
$ws.Range("I132").Value = 8717.710999999999
$ws.Range("K132").Value = 26153.133
$ws.Range("M132").Value = -23623.133

$ws.Range("H136").Value = 11889.917
$ws.Range("I136").Value = 8085.875
$ws.Range("J136").Value = 19498
$ws.Range("K136").Value = 24257.625
$ws.Range("L136").Value = 58494
$ws.Range("M136").Value = -21707.625
$ws.Range("N136").Value = -63594

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5439.4287
$ws.Range("I3").Value = 4035.1738
$ws.Range("K3").Value = 4035.1738
$ws.Range("M3").Value = -3921.1738

$ws.Range("H86").Value = 772315.6
$ws.Range("I86").Value = 1431463.9
$ws.Range("J86").Value = 3309.3333
$ws.Range("K86").Value = 1431463.9
$ws.Range("L86").Value = 3309.3333
$ws.Range("M86").Value = -1430340.9
$ws.Range("N86").Value = -5555.3333

$ws.Range("H89").Value = 772315.6
$ws.Range("I89").Value = 1431463.9
$ws.Range("J89").Value = 3309.3333
$ws.Range("K89").Value = 7157319.5
$ws.Range("L89").Value = 16546.6665
$ws.Range("M89").Value = -7151703.5
$ws.Range("N89").Value = -27778.6665

$ws.Range("H94").Value = 523.1622
$ws.Range("I94").Value = 472.39285
$ws.Range("K94").Value = 472.39285
$ws.Range("M94").Value = -21.39285000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3049.7
$ws.Range("J12").Value = 2671
$ws.Range("L12").Value = 2671
$ws.Range("N12").Value = -3011

$ws.Range("H15").Value = 5639.6
$ws.Range("J15").Value = 5639.6
$ws.Range("L15").Value = 5639.6
$ws.Range("N15").Value = -5979.6

$ws.Range("H16").Value = 1700.96
$ws.Range("I16").Value = 1671.45
$ws.Range("J16").Value = 1819
$ws.Range("K16").Value = 1671.45
$ws.Range("L16").Value = 1819
$ws.Range("M16").Value = -1384.45
$ws.Range("N16").Value = -2393

$ws.Range("H113").Value = 1700.96
$ws.Range("I113").Value = 1671.45
$ws.Range("J113").Value = 1819
$ws.Range("K113").Value = 1671.45
$ws.Range("L113").Value = 1819
$ws.Range("M113").Value = 498.55
$ws.Range("N113").Value = -6159

$ws.Range("H122").Value = 2154.182
$ws.Range("I122").Value = 2135.7693
$ws.Range("K122").Value = 6407.3079
$ws.Range("M122").Value = -3957.3079

$ws.Range("H132").Value = 15983.206
$ws.Range("I132").Value = 10673.692
$ws.Range("K132").Value = 32021.076
$ws.Range("M132").Value = -29491.076

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2269.3555
$ws.Range("I5").Value = 1510.9032
$ws.Range("J5").Value = 3948.7856
$ws.Range("K5").Value = 4532.7096
$ws.Range("L5").Value = 11846.3568
$ws.Range("M5").Value = -4420.7096
$ws.Range("N5").Value = -12070.3568

$ws.Range("H37").Value = 99968
$ws.Range("J37").Value = 99968
$ws.Range("L37").Value = 299904
$ws.Range("N37").Value = -300128

$ws.Range("H50").Value = 153.17392
$ws.Range("J50").Value = 150.2619
$ws.Range("L50").Value = 450.7857
$ws.Range("N50").Value = -1412.7857

$ws.Range("H53").Value = 153.17392
$ws.Range("J53").Value = 150.2619
$ws.Range("L53").Value = 450.7857
$ws.Range("N53").Value = -1412.7857

$ws.Range("H132").Value = 2076.7222
$ws.Range("I132").Value = 1891.4
$ws.Range("J132").Value = 2148
$ws.Range("K132").Value = 17022.6
$ws.Range("L132").Value = 19332
$ws.Range("M132").Value = -14492.6
$ws.Range("N132").Value = -24392

$ws.Range("H135").Value = 2269.3555
$ws.Range("I135").Value = 1510.9032
$ws.Range("J135").Value = 3948.7856
$ws.Range("K135").Value = 13598.1288
$ws.Range("L135").Value = 35539.0704
$ws.Range("M135").Value = -11063.1288
$ws.Range("N135").Value = -40609.0704

$ws.Range("H140").Value = 1753.0278
$ws.Range("I140").Value = 1639.7576
$ws.Range("K140").Value = 4919.2728
$ws.Range("M140").Value = 260.7272000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H80").Value = 11757.6
$ws.Range("I80").Value = 2096.3333
$ws.Range("K80").Value = 2096.3333
$ws.Range("M80").Value = -1098.3333

$ws.Range("H83").Value = 11757.6
$ws.Range("I83").Value = 2096.3333
$ws.Range("K83").Value = 10481.6665
$ws.Range("M83").Value = -5489.666499999999

$ws.Range("H97").Value = 867.19354
$ws.Range("I97").Value = 685
$ws.Range("K97").Value = 685
$ws.Range("M97").Value = -189

$ws.Range("H107").Value = 995.73334
$ws.Range("I107").Value = 995.1111
$ws.Range("J107").Value = 996.6667
$ws.Range("K107").Value = 995.1111
$ws.Range("L107").Value = 996.6667
$ws.Range("M107").Value = 924.8889
$ws.Range("N107").Value = -4836.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1336.2069
$ws.Range("I82").Value = 970.2381
$ws.Range("K82").Value = 970.2381
$ws.Range("M82").Value = -609.2381

$ws.Range("H85").Value = 1336.2069
$ws.Range("I85").Value = 970.2381
$ws.Range("K85").Value = 970.2381
$ws.Range("M85").Value = 277.7619

$ws.Range("H100").Value = 620752.8
$ws.Range("I100").Value = 3872
$ws.Range("J100").Value = 1854514.5
$ws.Range("K100").Value = 3872
$ws.Range("L100").Value = 1854514.5
$ws.Range("M100").Value = -3331
$ws.Range("N100").Value = -1855596.5

$ws.Range("H114").Value = 70000
$ws.Range("J114").Value = 70000
$ws.Range("L114").Value = 70000
$ws.Range("N114").Value = -78678

$ws.Range("H136").Value = 3178301
$ws.Range("I136").Value = 4834233.5
$ws.Range("K136").Value = 14502700.5
$ws.Range("M136").Value = -14500150.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 12
$ws.Range("I26").Value = 12
$ws.Range("K26").Value = 12
$ws.Range("M26").Value = 281

$ws.Range("H126").Value = 5970.1904
$ws.Range("I126").Value = 6324.9473
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 18974.8419
$ws.Range("L126").Value = 7800
$ws.Range("M126").Value = -16504.8419
$ws.Range("N126").Value = -12740

$ws.Range("H132").Value = 10844.116
$ws.Range("I132").Value = 8967.477000000001
$ws.Range("J132").Value = 30548.834
$ws.Range("K132").Value = 26902.431
$ws.Range("L132").Value = 91646.50199999999
$ws.Range("M132").Value = -24372.431
$ws.Range("N132").Value = -96706.50199999999
